$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H63").Value = 34545.168
$ws.Range("J63").Value = 34545.168
$ws.Range("L63").Value = 34545.168
$ws.Range("N63").Value = -35793.168
$ws.Range("H66").Value = 34545.168
$ws.Range("J66").Value = 34545.168
$ws.Range("L66").Value = 103635.504
$ws.Range("N66").Value = -109875.504
$ws.Range("H93").Value = 42500
$ws.Range("J93").Value = 42500
$ws.Range("L93").Value = 42500
$ws.Range("N93").Value = -47492
$ws.Range("H129").Value = 777.63495
$ws.Range("I129").Value = 311.55
$ws.Range("J129").Value = 994.4186
$ws.Range("K129").Value = 934.6500000000001
$ws.Range("L129").Value = 2983.2558
$ws.Range("M129").Value = 4065.35
$ws.Range("N129").Value = -12983.2558
$ws.Range("H130").Value = 35764.445
$ws.Range("J130").Value = 35764.445
$ws.Range("L130").Value = 35764.445
$ws.Range("N130").Value = -45804.445
$ws.Range("H137").Value = 1877.0741
$ws.Range("I137").Value = 1088.3889
$ws.Range("J137").Value = 3454.4443
$ws.Range("K137").Value = 3265.1667
$ws.Range("L137").Value = 10363.3329
$ws.Range("M137").Value = -715.1666999999998
$ws.Range("N137").Value = -15463.3329
$ws.Range("H138").Value = 2108.644
$ws.Range("I138").Value = 1783.091
$ws.Range("J138").Value = 2302.2163
$ws.Range("K138").Value = 5349.272999999999
$ws.Range("L138").Value = 6906.6489
$ws.Range("M138").Value = -209.2729999999992
$ws.Range("N138").Value = -17186.6489

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H55").Value = 15911.143
$ws.Range("J55").Value = 15911.143
$ws.Range("L55").Value = 15911.143
$ws.Range("N55").Value = -16541.143
$ws.Range("H61").Value = 1373.5
$ws.Range("I61").Value = 1309.2307
$ws.Range("K61").Value = 1309.2307
$ws.Range("M61").Value = -1097.2307
$ws.Range("H74").Value = 835.2258
$ws.Range("I74").Value = 758
$ws.Range("K74").Value = 758
$ws.Range("M74").Value = 116
$ws.Range("H77").Value = 835.2258
$ws.Range("I77").Value = 758
$ws.Range("K77").Value = 3790
$ws.Range("M77").Value = 578
$ws.Range("H80").Value = 21204
$ws.Range("J80").Value = 21204
$ws.Range("L80").Value = 21204
$ws.Range("N80").Value = -23200
$ws.Range("H83").Value = 21204
$ws.Range("J83").Value = 21204
$ws.Range("L83").Value = 63612
$ws.Range("N83").Value = -73596
$ws.Range("H103").Value = 36446
$ws.Range("J103").Value = 36446
$ws.Range("L103").Value = 36446
$ws.Range("N103").Value = -38790
$ws.Range("H113").Value = 53850
$ws.Range("J113").Value = 53850
$ws.Range("L113").Value = 53850
$ws.Range("N113").Value = -62528
$ws.Range("H123").Value = 56500
$ws.Range("J123").Value = 56500
$ws.Range("L123").Value = 56500
$ws.Range("N123").Value = -66300
$ws.Range("H129").Value = 48739.5
$ws.Range("J129").Value = 48739.5
$ws.Range("L129").Value = 48739.5
$ws.Range("N129").Value = -58739.5
$ws.Range("H131").Value = 39082.43
$ws.Range("J131").Value = 39082.43
$ws.Range("L131").Value = 39082.43
$ws.Range("N131").Value = -49162.43
$ws.Range("H134").Value = 37935.08
$ws.Range("J134").Value = 37935.08
$ws.Range("L134").Value = 37935.08
$ws.Range("N134").Value = -48075.08
$ws.Range("H135").Value = 40997
$ws.Range("J135").Value = 40997
$ws.Range("L135").Value = 40997
$ws.Range("N135").Value = -51137
$ws.Range("H136").Value = 1373.5
$ws.Range("I136").Value = 1309.2307
$ws.Range("K136").Value = 3927.6921
$ws.Range("M136").Value = -1377.6921

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H82").Value = 24921.572
$ws.Range("J82").Value = 30151.188
$ws.Range("L82").Value = 30151.188
$ws.Range("N82").Value = -30917.188
$ws.Range("H85").Value = 24921.572
$ws.Range("J85").Value = 30151.188
$ws.Range("L85").Value = 30151.188
$ws.Range("N85").Value = -32803.18799999999
$ws.Range("H122").Value = 29900
$ws.Range("J122").Value = 29900
$ws.Range("L122").Value = 29900
$ws.Range("N122").Value = -39700
$ws.Range("H135").Value = 59087.06
$ws.Range("J135").Value = 59087.06
$ws.Range("L135").Value = 59087.06
$ws.Range("N135").Value = -69227.06

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H122").Value = 1526.1428
$ws.Range("I122").Value = 1054
$ws.Range("J122").Value = 2255.818
$ws.Range("K122").Value = 3162
$ws.Range("L122").Value = 6767.454000000001
$ws.Range("M122").Value = -712
$ws.Range("N122").Value = -11667.454
$ws.Range("H127").Value = 52492.5
$ws.Range("J127").Value = 52492.5
$ws.Range("L127").Value = 52492.5
$ws.Range("N127").Value = -62412.5
$ws.Range("H141").Value = 46805.57
$ws.Range("J141").Value = 46805.57
$ws.Range("L141").Value = 46805.57
$ws.Range("N141").Value = -57165.57

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H92").Value = 796.1111
$ws.Range("I92").Value = 500
$ws.Range("J92").Value = 833.125
$ws.Range("K92").Value = 1500
$ws.Range("L92").Value = 2499.375
$ws.Range("M92").Value = -252
$ws.Range("N92").Value = -4995.375
$ws.Range("H113").Value = 790.6667
$ws.Range("I113").Value = 600
$ws.Range("K113").Value = 1800
$ws.Range("M113").Value = 370

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H43").Value = 4432.7144
$ws.Range("J43").Value = 9276.333000000001
$ws.Range("L43").Value = 9276.333000000001
$ws.Range("N43").Value = -9578.333000000001
$ws.Range("H127").Value = 48423
$ws.Range("J127").Value = 48423
$ws.Range("L127").Value = 48423
$ws.Range("N127").Value = -58343
$ws.Range("H128").Value = 53151.43
$ws.Range("J128").Value = 53151.43
$ws.Range("L128").Value = 53151.43
$ws.Range("N128").Value = -63111.43
$ws.Range("H132").Value = 2969
$ws.Range("I132").Value = 2395.9333
$ws.Range("J132").Value = 3685.3333
$ws.Range("K132").Value = 7187.7999
$ws.Range("L132").Value = 11055.9999
$ws.Range("M132").Value = -4657.7999
$ws.Range("N132").Value = -16115.9999
$ws.Range("H133").Value = 18838.334
$ws.Range("J133").Value = 18838.334
$ws.Range("L133").Value = 18838.334
$ws.Range("N133").Value = -28958.334
$ws.Range("H135").Value = 50200.715
$ws.Range("J135").Value = 50200.715
$ws.Range("L135").Value = 50200.715
$ws.Range("N135").Value = -60340.715

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H93").Value = 2306
$ws.Range("I93").Value = 2006.6666
$ws.Range("K93").Value = 2006.6666
$ws.Range("M93").Value = -758.6666
$ws.Range("H123").Value = 40325.4
$ws.Range("J123").Value = 40325.4
$ws.Range("L123").Value = 40325.4
$ws.Range("N123").Value = -50125.4
$ws.Range("H125").Value = 25000
$ws.Range("J125").Value = 25000
$ws.Range("L125").Value = 25000
$ws.Range("N125").Value = -34840
$ws.Range("H130").Value = 56340
$ws.Range("J130").Value = 56340
$ws.Range("L130").Value = 56340
$ws.Range("N130").Value = -66380
$ws.Range("H134").Value = 37119.332
$ws.Range("J134").Value = 37119.332
$ws.Range("L134").Value = 37119.332
$ws.Range("N134").Value = -47259.332

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H122").Value = 2129.577
$ws.Range("I122").Value = 1807.0416
$ws.Range("J122").Value = 6000
$ws.Range("K122").Value = 5421.1248
$ws.Range("L122").Value = 18000
$ws.Range("M122").Value = -2971.1248
$ws.Range("N122").Value = -22900
$ws.Range("H129").Value = 45429
$ws.Range("J129").Value = 45429
$ws.Range("L129").Value = 45429
$ws.Range("N129").Value = -55429
$ws.Range("H132").Value = 2141.2974
$ws.Range("I132").Value = 2221.1538
$ws.Range("J132").Value = 1952.5454
$ws.Range("K132").Value = 6663.4614
$ws.Range("L132").Value = 5857.6362
$ws.Range("M132").Value = -4133.4614
$ws.Range("N132").Value = -10917.6362
$ws.Range("H136").Value = 1198.9584
$ws.Range("I136").Value = 1015.2941
$ws.Range("J136").Value = 1645
$ws.Range("K136").Value = 3045.8823
$ws.Range("L136").Value = 4935
$ws.Range("M136").Value = -495.8822999999998
$ws.Range("N136").Value = -10035
